$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.511.43'
$ws.Range("E2").Value = '  -13.31%  '
$ws.Range("D3").Value = '2.321.91'
$ws.Range("E3").Value = '  -19.98%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '437.71'
$ws.Range("E5").Value = '  -16.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '121.21'
$ws.Range("E6").Value = '  -14.96%  '
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.468'
$ws.Range("E8").Value = '  -15.36%  '
$ws.Range("D9").Value = '2.313.72'
$ws.Range("E9").Value = '  -20.55%  '
$ws.Range("D10").Value = '5.21'
$ws.Range("E10").Value = '  -13.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0882'
$ws.Range("E11").Value = '  -18.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.303'
$ws.Range("E12").Value = '  -15.78%  '
$ws.Range("E13").Value = '  -5.41%  '
$ws.Range("D14").Value = '52.504.55'
$ws.Range("E14").Value = '  -13.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.88'
$ws.Range("E15").Value = '  -17.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000119'
$ws.Range("E16").Value = '  -15.66%  '
$ws.Range("D17").Value = '2.323.98'
$ws.Range("E17").Value = '  -20.21%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '302.32'
$ws.Range("E18").Value = '  -15.88%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.94'
$ws.Range("E19").Value = '  -21.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.93'
$ws.Range("E20").Value = '  -23.44%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.60'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.16'
$ws.Range("E23").Value = '  -22.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '53.53'
$ws.Range("E24").Value = '  -17.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.364'
$ws.Range("E25").Value = '  -19.76%  '
$ws.Range("D26").Value = '0.145'
$ws.Range("E26").Value = '  -20.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.94'
$ws.Range("E27").Value = '  -11.99%  '
$ws.Range("D28").Value = '0.995'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").Value = '0.0₃0675'
$ws.Range("E29").Value = '  -19.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '143.48'
$ws.Range("E30").Value = '  -4.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '16.99'
$ws.Range("E31").Value = '  -14.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.33'
$ws.Range("E32").Value = '  -20.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.73'
$ws.Range("E33").Value = '  -15.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.49'
$ws.Range("E34").Value = '  -19.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.821'
$ws.Range("E35").Value = '  -18.23%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.997'
$ws.Range("E36").Value = '  -17.01%  '
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.992'
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.97'
$ws.Range("E38").Value = '  -15.55%  '
$ws.Range("D39").Value = '10.16'
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.15'
$ws.Range("E40").Value = '  -15.25%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0501'
$ws.Range("E41").Value = '  -14.35%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.22'
$ws.Range("E42").Value = '  -18.00%  '
$ws.Range("D43").Value = '1.900.49'
$ws.Range("E43").Value = '  -17.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.524'
$ws.Range("E44").Value = '  -19.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0206'
$ws.Range("E45").Value = '  -13.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0824'
$ws.Range("E46").Value = '  -10.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.61'
$ws.Range("E47").Value = '  -23.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.94'
$ws.Range("E48").Value = '  -20.40%  '
$ws.Range("E49").Value = '  -5.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '15.11'
$ws.Range("E50").Value = '  -17.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.43'
$ws.Range("E51").Value = '  -14.50%  '
